$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.2042
$ws.Range("E4").Value = 16.4343
$ws.Range("B7").Value = 5.724
$ws.Range("A8").Value = -22.31150000000001
$ws.Range("A10").Value = -21.61899999999999
$ws.Range("E11").Value = 16.69129999999999
$ws.Range("A12").Value = -21.58219999999999
$ws.Range("B14").Value = 6.359000000000005
$ws.Range("E14").Value = 16.7216
$ws.Range("B15").Value = 4.895899999999997
$ws.Range("A18").Value = -21.76960000000001
$ws.Range("B18").Value = 6.068699999999995
$ws.Range("E18").Value = 17.88860000000003
$ws.Range("E19").Value = 16.4435
$ws.Range("B20").Value = 8.836200000000005
$ws.Range("E21").Value = 16.53810000000001
$ws.Range("A25").Value = -21.50889999999999
$ws.Range("E27").Value = 16.47699999999999
$ws.Range("B29").Value = 5.195300000000004
$ws.Range("B30").Value = 5.661500000000002
$ws.Range("B31").Value = 5.769400000000005
$ws.Range("E31").Value = 16.62829999999998
$ws.Range("B35").Value = 8.627400000000007
$ws.Range("A37").Value = -19.17989999999999
$ws.Range("E38").Value = 16.21909999999999
$ws.Range("B40").Value = 9.248499999999989
$ws.Range("E42").Value = 16.35360000000001
$ws.Range("B44").Value = 4.531600000000003
$ws.Range("E44").Value = 16.5221
$ws.Range("E47").Value = 16.5156
$ws.Range("B50").Value = 4.395800000000001
$ws.Range("B54").Value = 4.105599999999999
$ws.Range("A55").Value = -21.9119
$ws.Range("E56").Value = 16.42840000000001
$ws.Range("E58").Value = 16.16950000000001
$ws.Range("E65").Value = 17.32310000000001
$ws.Range("A68").Value = -21.48669999999999
$ws.Range("B68").Value = 4.532099999999999
$ws.Range("E73").Value = 17.46910000000001
$ws.Range("B76").Value = 5.790900000000001
$ws.Range("A77").Value = -20.12539999999998
$ws.Range("A78").Value = -20.33309999999997
$ws.Range("A79").Value = -20.07339999999999
$ws.Range("A80").Value = -19.68509999999999
$ws.Range("A81").Value = -21.7786
$ws.Range("A82").Value = -21.7832
$ws.Range("A84").Value = -22.0121
$ws.Range("B87").Value = 4.647199999999998
$ws.Range("B88").Value = 4.529799999999997
$ws.Range("E90").Value = 16.3458
$ws.Range("B92").Value = 4.846099999999997
$ws.Range("E92").Value = 18.62060000000001
$ws.Range("E94").Value = 18.95300000000001
$ws.Range("E95").Value = 17.99350000000002
$ws.Range("B96").Value = 5.216200000000008
$ws.Range("B98").Value = 5.839499999999999
$ws.Range("A101").Value = -21.98140000000001
$ws.Range("B101").Value = 5.937800000000002
$ws.Range("E101").Value = 16.9359
$ws.Range("A102").Value = -21.6338
$ws.Range("B102").Value = 5.628900000000004
